$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 72377.28999999999
$ws.Range("J17").Value = 75047.55499999999
$ws.Range("L17").Value = 225142.665
$ws.Range("N17").Value = -225478.665
$ws.Range("H53").Value = 14493580
$ws.Range("I53").Value = 66667776
$ws.Range("K53").Value = 66667776
$ws.Range("M53").Value = -66667139
$ws.Range("H57").Value = 51992.5
$ws.Range("J57").Value = 51992.5
$ws.Range("L57").Value = 155977.5
$ws.Range("N57").Value = -156975.5
$ws.Range("H138").Value = 7104.9688
$ws.Range("I138").Value = 8098.5
$ws.Range("J138").Value = 7038.7334
$ws.Range("K138").Value = 24295.5
$ws.Range("L138").Value = 21116.2002
$ws.Range("M138").Value = -19155.5
$ws.Range("N138").Value = -31396.2002

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4463
$ws.Range("I2").Value = 3579.7
$ws.Range("K2").Value = 3579.7
$ws.Range("M2").Value = -3466.7
$ws.Range("H32").Value = 1927.0793
$ws.Range("I32").Value = 1981.1803
$ws.Range("J32").Value = 277
$ws.Range("K32").Value = 1981.1803
$ws.Range("L32").Value = 277
$ws.Range("M32").Value = -1694.1803
$ws.Range("N32").Value = -851
$ws.Range("H61").Value = 8815.6
$ws.Range("I61").Value = 8818.666999999999
$ws.Range("J61").Value = 8814.286
$ws.Range("K61").Value = 8818.666999999999
$ws.Range("L61").Value = 8814.286
$ws.Range("M61").Value = -8606.666999999999
$ws.Range("N61").Value = -9238.286
$ws.Range("H74").Value = 12474.259
$ws.Range("I74").Value = 14984.7
$ws.Range("J74").Value = 5301.5713
$ws.Range("K74").Value = 14984.7
$ws.Range("L74").Value = 5301.5713
$ws.Range("M74").Value = -14110.7
$ws.Range("N74").Value = -7049.5713
$ws.Range("H77").Value = 12474.259
$ws.Range("I77").Value = 14984.7
$ws.Range("J77").Value = 5301.5713
$ws.Range("K77").Value = 74923.5
$ws.Range("L77").Value = 26507.8565
$ws.Range("M77").Value = -70555.5
$ws.Range("N77").Value = -35243.85649999999
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H116").Value = 4463
$ws.Range("I116").Value = 3579.7
$ws.Range("K116").Value = 3579.7
$ws.Range("M116").Value = -1285.7
$ws.Range("H122").Value = 4944.4287
$ws.Range("I122").Value = 4449.737
$ws.Range("K122").Value = 13349.211
$ws.Range("M122").Value = -10899.211
$ws.Range("H136").Value = 8815.6
$ws.Range("I136").Value = 8818.666999999999
$ws.Range("J136").Value = 8814.286
$ws.Range("K136").Value = 26456.001
$ws.Range("L136").Value = 26442.858
$ws.Range("M136").Value = -23906.001
$ws.Range("N136").Value = -31542.858

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4463
$ws.Range("I3").Value = 3579.7
$ws.Range("K3").Value = 3579.7
$ws.Range("M3").Value = -3465.7
$ws.Range("H134").Value = 49205.875
$ws.Range("I134").Value = 8288.857
$ws.Range("J134").Value = 106489.7
$ws.Range("K134").Value = 24866.571
$ws.Range("L134").Value = 319469.1
$ws.Range("M134").Value = -22331.571
$ws.Range("N134").Value = -324539.1

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4009.375
$ws.Range("I31").Value = 1725
$ws.Range("K31").Value = 1725
$ws.Range("M31").Value = -1430
$ws.Range("H34").Value = 4009.375
$ws.Range("I34").Value = 1725
$ws.Range("K34").Value = 1725
$ws.Range("M34").Value = -1523
$ws.Range("H58").Value = 913868.4
$ws.Range("I58").Value = 1431793.1
$ws.Range("K58").Value = 1431793.1
$ws.Range("M58").Value = -1431590.1
$ws.Range("H134").Value = 592628.25
$ws.Range("I134").Value = 4152.727
$ws.Range("K134").Value = 12458.181
$ws.Range("M134").Value = -9923.181
$ws.Range("H136").Value = 913868.4
$ws.Range("I136").Value = 1431793.1
$ws.Range("K136").Value = 4295379.300000001
$ws.Range("M136").Value = -4292829.300000001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 128.16667
$ws.Range("I23").Value = 150.8
$ws.Range("J23").Value = 15
$ws.Range("K23").Value = 452.4
$ws.Range("L23").Value = 45
$ws.Range("M23").Value = -217.4
$ws.Range("N23").Value = -515
$ws.Range("H68").Value = 1354.2778
$ws.Range("J68").Value = 1598.3077
$ws.Range("L68").Value = 4794.9231
$ws.Range("N68").Value = -6416.9231
$ws.Range("H71").Value = 1354.2778
$ws.Range("J71").Value = 1598.3077
$ws.Range("L71").Value = 14384.7693
$ws.Range("N71").Value = -22496.7693
$ws.Range("H113").Value = 6173673
$ws.Range("J113").Value = 1000.2
$ws.Range("L113").Value = 3000.6
$ws.Range("N113").Value = -7340.6
$ws.Range("H132").Value = 3362.3635
$ws.Range("I132").Value = 2082.6667
$ws.Range("K132").Value = 18744.0003
$ws.Range("M132").Value = -16214.0003
$ws.Range("H140").Value = 3153.8823
$ws.Range("I140").Value = 3153.8823
$ws.Range("K140").Value = 9461.6469
$ws.Range("M140").Value = -4281.6469

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 1811.8
$ws.Range("I99").Value = 1811.8
$ws.Range("K99").Value = 1811.8
$ws.Range("M99").Value = 434.2
$ws.Range("H132").Value = 284635.28
$ws.Range("J132").Value = 36950.414
$ws.Range("L132").Value = 110851.242
$ws.Range("N132").Value = -115911.242

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 81860.87
$ws.Range("I100").Value = 133939.12
$ws.Range("K100").Value = 133939.12
$ws.Range("M100").Value = -133398.12
$ws.Range("H132").Value = 6433.5557
$ws.Range("I132").Value = 4999.6665
$ws.Range("J132").Value = 7150.5
$ws.Range("K132").Value = 14998.9995
$ws.Range("L132").Value = 21451.5
$ws.Range("M132").Value = -12468.9995
$ws.Range("N132").Value = -26511.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 956.1429000000001
$ws.Range("I113").Value = 843.5833
$ws.Range("K113").Value = 2530.7499
$ws.Range("M113").Value = -360.7498999999998
$ws.Range("H132").Value = 38346.9
$ws.Range("I132").Value = 2839.3157
$ws.Range("K132").Value = 8517.947100000001
$ws.Range("M132").Value = -5987.947100000001
$ws.Range("H133").Value = 65426.25
$ws.Range("J133").Value = 65426.25
$ws.Range("L133").Value = 65426.25
$ws.Range("N133").Value = -75546.25
